$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Yerba Buena, 30 de Octubre de 1986" paragraph
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# ---------------------------------------------------------------------
# 2) "ORDENANZA N\u00ba 211" paragraph
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# 3) Split "VISTO: ..." into its own label paragraph + body paragraph
# ---------------------------------------------------------------------
$visto = $d.Range(0, 0)
$vistoFind = $visto.Find
$vistoFind.Execute("VISTO: ")
$visto.Collapse(0)
$visto.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 12
$p3.Format.SpaceAfter = 6
$p3.Format.Alignment = 0
$p3.Range.Font.Bold = $true

$p4 = $d.Paragraphs.Item(4)
$p4.Format.KeepWithNext = $true
$p4.Format.SpaceAfter = 6
$p4.Format.Alignment = 0
$p4Start = $p4.Range.Duplicate
$p4Start.Collapse(1)
$p4Start.InsertBefore(" ")

# ---------------------------------------------------------------------
# 4) Split "CONSIDERANDO: ..." into its own label paragraph + body paragraph
# ---------------------------------------------------------------------
$consid = $d.Range(0, 0)
$considFind = $consid.Find
$considFind.Execute("CONSIDERANDO: ")
$consid.Collapse(0)
$consid.InsertParagraphAfter()

$p5 = $d.Paragraphs.Item(5)
$p5.Format.KeepWithNext = $true
$p5.Format.SpaceBefore = 12
$p5.Format.SpaceAfter = 6
$p5.Format.Alignment = 0
$p5.Range.Font.Bold = $true

$p6 = $d.Paragraphs.Item(6)
$p6.Format.KeepWithNext = $true
$p6.Format.SpaceAfter = 6
$p6.Format.Alignment = 0
$p6Start = $p6.Range.Duplicate
$p6Start.Collapse(1)
$p6Start.InsertBefore(" ")

# ---------------------------------------------------------------------
# 5) The six tabbed "Que ..." paragraphs (now items 7..11) + "Por ello:" (12)
# ---------------------------------------------------------------------
for ($i = 7; $i -le 12; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Format.KeepWithNext = $true
    $p.Format.SpaceAfter = 6
    $p.Format.Alignment = 0
}

# ---------------------------------------------------------------------
# 6) "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA" (13)
# ---------------------------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$p13.Format.KeepWithNext = $true
$p13.Format.SpaceBefore = 18
$p13.Format.SpaceAfter = 18
$p13.Format.LeftIndent = 99.2
$p13.Format.RightIndent = 99.2
$p13.Range.Font.Bold = $true

# ---------------------------------------------------------------------
# 7) ARTICULO PRIMERO / SEGUNDO / TERCERO / CUARTO (14..17)
# ---------------------------------------------------------------------
for ($i = 14; $i -le 17; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Format.KeepWithNext = $true
    $p.Format.SpaceAfter = 6
    $p.Format.Alignment = 0
}

$articuloRange = $d.Range(0, 0)
$articuloFind = $articuloRange.Find
$labels = @("ARTICULO PRIMERO:", "ARTICULO SEGUNDO:", "ARTICULO TERCERO:", "ARTICULO CUARTO:")
foreach ($label in $labels) {
    $articuloFind.Execute($label)
    $articuloRange.Font.Underline = 1
    $articuloRange.Collapse(0)
}

# ---------------------------------------------------------------------
# 8) Section: default footer + page numbering restart at 127
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$ftr.PageNumbers.StartingNumber = 127

$ftrXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Piedepgina"/><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:color w:val="808080" w:themeColor="background1" w:themeShade="80"/><w:sz w:val="20"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ftr.Range.InsertXML($ftrXml)

Write-Output "edit complete"
